$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Update "Wales_Bowel_Uptake" sheet: append two more years of data
# ------------------------------------------------------------------
$waleBowel = $wb.Worksheets.Item("Wales_Bowel_Uptake")
$waleBowel.Range("A10").Value = 2020
$waleBowel.Range("B10").Value = 67.1
$waleBowel.Range("A11").Value = 2021
$waleBowel.Range("B11").Value = 67.2
$waleBowel.Range("B12").Select() | Out-Null

# ------------------------------------------------------------------
# 2. Update "NIreland_Breast_Uptake" sheet: change the active selection
#    to the header row A1:B1 (no data changes here)
# ------------------------------------------------------------------
$niBreast = $wb.Worksheets.Item("NIreland_Breast_Uptake")
$niBreast.Range("A1:B1").Select() | Out-Null

# ------------------------------------------------------------------
# 3. Add a brand new sheet "NIreland_Bowel_Uptake" at the end of the
#    workbook with the bowel screening uptake figures for NI
# ------------------------------------------------------------------
$niBowel = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$niBowel.Name = "NIreland_Bowel_Uptake"

$niBowel.Range("A1").Value = "Year"
$niBowel.Range("B1").Value = "Uptake"

$niBowel.Range("A2").Value = 2018
$niBowel.Range("B2").Value = 54.3

$niBowel.Range("A3").Value = 2019
$niBowel.Range("B3").Value = 53.04

$niBowel.Range("A4").Value = 2020
$niBowel.Range("B4").Value = 57.94

$niBowel.Range("A5").Value = 2021
$niBowel.Range("B5").Value = 62.15

$niBowel.Range("A6").Value = 2022
$niBowel.Range("B6").Value = 57.85

# Make this newly added sheet the active one, with the cursor resting
# on K12 (matching the recorded cursor position) and make it the
# visible/selected tab.
$niBowel.Activate() | Out-Null
$niBowel.Range("K12").Select() | Out-Null
